$wb = $excel.ActiveWorkbook

$wsWith    = $wb.Worksheets.Item("pages_with_related_resources")
$wsWithout = $wb.Worksheets.Item("pages_without_related_resources")
$wsExt     = $wb.Worksheets.Item("pages_with_external_resources")

# --- Data fix: corrected Spanish URL path (affected by site migration) ---
$oldPath = "espanol/news-events/cancer-currents-blog/2019/vitamina-d-complemento-cancer-prevencion"
$newPath = "espanol/noticias/temas-y-relatos-blog/2019/vitamina-d-complemento-cancer-prevencion"

$wsWith.Range("A4").Value = $newPath
$wsExt.Range("A9").Value = $newPath

# --- Remove the stray/no-op alignment formatting left on A3:A5 of the
#     "pages_without_related_resources" sheet (applyAlignment with no actual
#     alignment set - a leftover that never rendered any different from the
#     default style) ---
$wsWithout.Range("A3:A5").ClearFormats()

# --- Selection bookkeeping: each sheet had been fully selected (Ctrl+A) while
#     reviewing the corrected rows ---
$wsWith.Activate()
$wsWith.Range("A1:XFD8").Select()

$wsWithout.Activate()
$wsWithout.Range("A1:XFD5").Select()

$wsExt.Activate()
$wsExt.Range("A1:XFD9").Select()

Write-Host "Applied related-resources URL fix and formatting cleanup."
